$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.336.32"
$ws.Range("E2").Value = "  -3.97%  "
$ws.Range("D3").Value = "3.825.22"
$ws.Range("E3").Value = "  -4.49%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'" + "590.39"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'" + "166.11"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("D7").Value = "'" + "0.661"
$ws.Range("E7").Value = "  -3.47%  "
$ws.Range("D8").Value = "'" + "1.00"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'" + "0.738"
$ws.Range("E9").Value = "  -1.80%  "
$ws.Range("D10").Value = "'" + "0.173"
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("D11").Value = "'" + "52.46"
$ws.Range("E11").Value = "  -3.07%  "
$ws.Range("D12").Value = "'" + "0.0000317"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "'" + "11.21"
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("D14").Value = "4.428.26"
$ws.Range("E14").Value = "  -4.57%  "
$ws.Range("D15").Value = "3.825.05"
$ws.Range("E15").Value = "  -4.56%  "
$ws.Range("D16").Value = "'" + "20.61"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").Value = "'" + "13.70"
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("E18").Value = "  -6.44%  "
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").Value = "70.172.45"
$ws.Range("E20").Value = "  -3.79%  "
$ws.Range("D21").Value = "'" + "432.51"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "'" + "4.72"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "'" + "92.93"
$ws.Range("E23").Value = "  -3.57%  "
$ws.Range("D24").Value = "'" + "3.26"
$ws.Range("E24").Value = "  -5.09%  "
$ws.Range("D25").Value = "'" + "13.73"
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("D26").Value = "'" + "11.44"
$ws.Range("E26").Value = "  +1.77%  "
$ws.Range("D27").Value = "'" + "3.96"
$ws.Range("E27").Value = "  -10.32%  "
$ws.Range("D28").Value = "'" + "5.96"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'" + "10.41"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").Value = "'" + "34.71"
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("D31").Value = "'" + "8.13"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("D32").Value = "'" + "13.37"
$ws.Range("E32").Value = "  -3.30%  "
$ws.Range("D33").Value = "'" + "47.52"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("D35").Value = "'" + "0.0" + [char]0x2083 + "0986"
$ws.Range("E35").Value = "  +8.08%  "
$ws.Range("D36").Value = "'" + "68.31"
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("D37").Value = "'" + "638.87"
$ws.Range("E37").Value = "  -4.83%  "
$ws.Range("D38").Value = "'" + "0.428"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").Value = "'" + "0.145"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'" + "3.30"
$ws.Range("E41").Value = "  +25.56%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'" + "0.999"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "'" + "2.93"
$ws.Range("E43").Value = "  +11.84%  "
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("D45").Value = "'" + "0.0465"
$ws.Range("E45").Value = "  -5.77%  "
$ws.Range("D46").Value = "'" + "9.82"
$ws.Range("E46").Value = "  -7.93%  "
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("D48").Value = "2.855.70"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "'" + "2.78"
$ws.Range("E49").Value = "  -16.87%  "
$ws.Range("D50").Value = "'" + "3.23"
$ws.Range("E50").Value = "  -5.20%  "
$ws.Range("D51").Value = "'" + "0.000273"
$ws.Range("E51").Value = "  +0.11%  "
